# Populate the "Personal Data" sheet with rows for the two supervisors and
# their four direct reports (data rows 3-8, below the existing title/header
# rows 1-2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "dd/mm/yyyy"

function Set-PersonRow($Row, $Name, $Id, $Age, $JoinDate, $IsSupervisor) {
    # A: name (first_name column holds the full placeholder name)
    $ws.Cells.Item($Row, 1).Value = $Name
    # B: id
    $ws.Cells.Item($Row, 2).Value = $Id
    # C: last_name placeholder
    $ws.Cells.Item($Row, 3).Value = "__"
    # D: age
    $ws.Cells.Item($Row, 4).Value = $Age
    # E: join_date (value first, number format applied after the name-cell
    # styling below so the style table is built up in the same order as the
    # source workbook: name-alignment style, then date style, then indent).
    $ws.Cells.Item($Row, 5).Value = $JoinDate

    $nameCell = $ws.Cells.Item($Row, 1)
    $nameCell.Font.Bold = $true
    $nameCell.WrapText = $true
    $nameCell.HorizontalAlignment = -4131

    if ($IsSupervisor) {
        $nameCell.VerticalAlignment = -4108
        $ws.Rows.Item($Row).RowHeight = 33.75
    } else {
        $nameCell.IndentLevel = 4
    }

    $ws.Cells.Item($Row, 5).NumberFormat = $dateFormat
}

function Set-YearRatings($Row, $SupervisorCol, $ClientsCol, $AiCol, $DateCol, $SupervisorRating, $ClientsRating, $AiRating, $RatingDate) {
    $ws.Range($SupervisorCol + $Row).Value = $SupervisorRating
    $ws.Range($ClientsCol + $Row).Value = $ClientsRating
    $ws.Range($AiCol + $Row).Value = $AiRating
    $ws.Range($DateCol + $Row).Value = $RatingDate
    $ws.Range($DateCol + $Row).NumberFormat = $dateFormat
}

# Row 3: Supervisor_2 (no historical ratings on file)
Set-PersonRow 3 "Supervisor_2" 10 20 44197 $true

# Row 4: employee_3, rated in 2020 (columns N:Q)
Set-PersonRow 4 "employee_3" 4 23 40909 $false
Set-YearRatings 4 "N" "O" "P" "Q" 7.3 4 8 43831

# Row 5: employee_4, rated in 2020 (columns N:Q)
Set-PersonRow 5 "employee_4" 5 23 40909 $false
Set-YearRatings 5 "N" "O" "P" "Q" 2.5 2.4 10 43831

# Row 6: Supervisor_1, rated in 2022 (columns F:I)
Set-PersonRow 6 "Supervisor_1" 1 20 44197 $true
Set-YearRatings 6 "F" "G" "H" "I" 9.5 8.6 8 44562

# Row 7: employee_1, rated in 2022 (columns F:I)
Set-PersonRow 7 "employee_1" 2 22 44562 $false
Set-YearRatings 7 "F" "G" "H" "I" 8.3 8.6 8 44562

# Row 8: employee_2, rated in 2021 (columns J:M)
Set-PersonRow 8 "employee_2" 3 12 40909 $false
Set-YearRatings 8 "J" "K" "L" "M" 8.3 1 8 44197
